# "json schema of mock user" -- adds a new mock-user record (phone, name,
# email, dob, card number, IDEP id) and repoints every cell that referenced
# the previous mock-user record's fields to the new ones. Also bumps the
# numeric id in RegisterUser!M2.

$wb = $excel.ActiveWorkbook

$newPhone = "6818087119"
$newName  = "Niyati"
$newEmail = "KailashRaja65125@example.net"
$newDob   = "1989-01-26"
$newCard  = "752745609345715"
$newIdep  = "IDEP6169959956125LUT"

# Helper: write a string value into a cell *as text*, even when the string
# looks like a number (plain .Value assignment would otherwise silently
# coerce an all-digit string like a phone/card number into a numeric
# cell). Round-tripping it through a text formula and then pasting back
# as a value keeps the cell's existing style/number format untouched.
function Set-TextValue($range, $text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# SendOtp!C2
$ws = $wb.Worksheets.Item("SendOtp")
Set-TextValue $ws.Range("C2") $newPhone

# User_Authenticate!C2:C6
$ws = $wb.Worksheets.Item("User_Authenticate")
Set-TextValue $ws.Range("C2") $newPhone
Set-TextValue $ws.Range("C3") $newPhone
Set-TextValue $ws.Range("C4") $newPhone
Set-TextValue $ws.Range("C5") $newPhone
Set-TextValue $ws.Range("C6") $newPhone

# UserOnboarding!E2
$ws = $wb.Worksheets.Item("UserOnboarding")
Set-TextValue $ws.Range("E2") $newCard

# Create_Bnpl_Transaction!G2
$ws = $wb.Worksheets.Item("Create_Bnpl_Transaction")
Set-TextValue $ws.Range("G2") $newCard

# UpdateUser!C2
$ws = $wb.Worksheets.Item("UpdateUser")
$ws.Range("C2").Value = $newName

# " UpdateUser"!E2:G2
$ws = $wb.Worksheets.Item(" UpdateUser")
$ws.Range("E2").Value = $newName
$ws.Range("F2").Value = $newEmail
$ws.Range("G2").Value = $newDob

# RegisterUser!D2,E2,I2,L2,M2
$ws = $wb.Worksheets.Item("RegisterUser")
Set-TextValue $ws.Range("D2") $newPhone
$ws.Range("E2").Value = $newEmail
Set-TextValue $ws.Range("I2") $newCard
$ws.Range("L2").Value = $newIdep
$ws.Range("M2").Value = 4202229.0
